$d = $word.ActiveDocument

# Replace header cell text to reflect the new naming convention.
# wdFindContinue = 1, wdReplaceAll = 2

$d.Content.Find.Execute("Case Sample Name", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Treated Sample Name", 2)

$d.Content.Find.Execute("Case Chromosome", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Chromosome", 2)

$d.Content.Find.Execute("Case Event Start", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Event Start", 2)

$d.Content.Find.Execute("Case Event End", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Event End", 2)

$d.Content.Find.Execute("Case Event Size", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Event Size", 2)

$d.Content.Find.Execute("Case Fractional Copy Number", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Treated Fractional Copy Number", 2)
